$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells whose new values look numeric stay as plain text (matches source formatting)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values from the crypto data refresh
$ws.Range("D2").Value = "65.231.95"
$ws.Range("E2").Value = "  +2.84%  "
$ws.Range("D3").Value = "2.653.40"
$ws.Range("E3").Value = "  +2.64%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "596.18"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("D6").Value = "156.54"
$ws.Range("E6").Value = "  +3.72%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +1.13%  "
$ws.Range("E9").Value = "  +7.55%  "
$ws.Range("E10").Value = "  +3.89%  "
$ws.Range("E11").Value = "  +1.58%  "
$ws.Range("E12").Value = "  +1.82%  "
$ws.Range("D13").Value = "29.06"
$ws.Range("E13").Value = "  +5.26%  "
$ws.Range("D14").Value = "0.0000185"
$ws.Range("E14").Value = "  +18.45%  "
$ws.Range("D15").Value = "3.130.89"
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("D16").Value = "65.109.96"
$ws.Range("E16").Value = "  +3.03%  "
$ws.Range("D17").Value = "2.656.81"
$ws.Range("E17").Value = "  +3.15%  "
$ws.Range("D18").Value = "12.63"
$ws.Range("E18").Value = "  +3.00%  "
$ws.Range("D19").Value = "4.81"
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("D20").Value = "355.22"
$ws.Range("E20").Value = "  +2.80%  "
$ws.Range("D21").Value = "7.29"
$ws.Range("E21").Value = "  +6.04%  "
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").Value = "68.27"
$ws.Range("E23").Value = "  +1.21%  "
$ws.Range("D24").Value = "1.70"
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").Value = "9.52"
$ws.Range("E25").Value = "  +2.87%  "
$ws.Range("E26").Value = "  -1.84%  "
$ws.Range("E27").Value = "  +1.09%  "
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("E29").Value = "  +11.41%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "520.93"
$ws.Range("E31").Value = "  -7.58%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "2.10"
$ws.Range("E32").Value = "  +3.30%  "
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("D34").Value = "5.63"
$ws.Range("E34").Value = "  +7.14%  "
$ws.Range("D35").Value = "6.37"
$ws.Range("E35").Value = "  +3.93%  "
$ws.Range("D36").Value = "0.429"
$ws.Range("E36").Value = "  +3.81%  "
$ws.Range("D37").Value = "164.96"
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "20.33"
$ws.Range("E38").Value = "  +3.93%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "2.03"
$ws.Range("E39").Value = "  +4.98%  "
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").Value = "42.16"
$ws.Range("E42").Value = "  +6.55%  "
$ws.Range("D43").Value = "165.83"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").Value = "4.11"
$ws.Range("E44").Value = "  +3.00%  "
$ws.Range("D45").Value = "0.0619"
$ws.Range("E45").Value = "  +5.72%  "
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("E47").Value = "  +5.04%  "
$ws.Range("E48").Value = "  +3.40%  "
$ws.Range("E49").Value = "  +1.61%  "
$ws.Range("D50").Value = "0.0987"
$ws.Range("E50").Value = "  +2.69%  "
$ws.Range("D51").Value = "19.54"
$ws.Range("E51").Value = "  +1.80%  "

Write-Output "Applied 91 cell updates"
